# Commit: "include Liyunet comments on Oromiya"
#
# This script reproduces, in terms of observable cell values, the effect of
# the OOXML diff against xl/sharedStrings.xml and xl/worksheets/sheet1.xml:
#   - "missing"    -> "meskan"          (affects the Wereda column wherever it
#                                         was used, incl. rows that used to
#                                         read the separate "misrak meskan")
#   - "halaba"     -> "siltie"          (Zone column, rows 281-328)
#   - "wera"       -> "alaba special"   (Wereda column, rows 281-328, and it
#                                         also absorbs the rows that used to
#                                         read "atoti hullo" / "wera dijjo")
#   - "ho/kuke"        -> "holugeb kuke"
#   - "ta/bedene"      -> "tachegnawo bedane"
#   - "la/bedene"      -> "laygnawo bedane"
#   - "mekala ha"      -> "huletegna mekala"
#   - "la/lenda"       -> "layegnawo lenda"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "missing" -> "meskan" -------------------------------------------------
# Rows 2-14: Wereda column (C) already reads "missing".
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("C$r").Value = "meskan"
}

# Rows 84-98: Wereda column (C) used to read the distinct string
# "misrak meskan"; it now points at the same "meskan" string as above.
for ($r = 84; $r -le 98; $r++) {
    $ws.Range("C$r").Value = "meskan"
}

# --- Halaba zone relabelled as Siltie / Alaba special (rows 281-328) ------
for ($r = 281; $r -le 328; $r++) {
    $ws.Range("B$r").Value = "siltie"
    $ws.Range("C$r").Value = "alaba special"
}

# --- A handful of Kebele-name (D column) spelling/content fixes -----------
$ws.Range("D285").Value = "holugeb kuke"
$ws.Range("D287").Value = "tachegnawo bedane"
$ws.Range("D289").Value = "laygnawo bedane"
$ws.Range("D299").Value = "huletegna mekala"
$ws.Range("D303").Value = "layegnawo lenda"
